$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new rows of trial data (T8, T9, T10), reusing the existing
# shared-string entries and matching the formatting of the preceding rows.
$ws.Range("A15").Value = "T8"
$ws.Range("B15").Value = 700
$ws.Range("C15").Value = 12
$ws.Range("D15").Value = 0.58
$ws.Range("E15").Value = 0.15
$ws.Range("F15").Value = 0.23
$ws.Range("G15").Value = 0.29
$ws.Range("H15").Value = 0.15

$ws.Range("A16").Value = "T9"
$ws.Range("B16").Value = 660
$ws.Range("C16").Value = 14
$ws.Range("D16").Value = 0.47
$ws.Range("E16").Value = 0.2
$ws.Range("F16").Value = 0.11
$ws.Range("G16").Value = 0.26
$ws.Range("H16").Value = 0.39

$ws.Range("A17").Value = "T10"
$ws.Range("B17").Value = 760
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 0.42
$ws.Range("E17").Value = 0.17
$ws.Range("F17").Value = 0.21
$ws.Range("G17").Value = 0.23
$ws.Range("H17").Value = 0.29

# Apply the same formatting used by the rest of the trial-data block
# (rows 8-14) to the newly-added rows.
$ws.Range("A14:H14").Copy()
$ws.Range("A15:H17").PasteSpecial(-4122)  # xlPasteFormats

# Extend column P with the same (blank, formatted) cells as column O for
# the first four header/summary rows.
$ws.Range("O1:O4").Copy()
$ws.Range("P1:P4").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Update the active selection to reflect the new end-of-data cell.
$ws.Range("A18").Select()
